# Uber Dispatching Data - "Added Uber minicase data"
# Flip the three match-indicator flags from 0 -> 1 on Josephines_Match
# and mirror the same rider flags on Riders, then leave the UI selection
# state (active sheet / selected cell per sheet) the way the saved
# workbook had it: Drivers active (instead of Josephines_Match), with
# new last-selected cells on each sheet.

$wb = $excel.ActiveWorkbook

# --- Josephines_Match: set three binary match flags to 1 ---
$wsMatch = $wb.Worksheets.Item("Josephines_Match")
[void]$wsMatch.Activate()
$wsMatch.Range("C9").Value = 1
$wsMatch.Range("H10").Value = 1
$wsMatch.Range("B11").Value = 1
[void]$wsMatch.Range("P15").Select()

# --- Riders: set the corresponding match flags to 1 ---
$wsRiders = $wb.Worksheets.Item("Riders")
[void]$wsRiders.Activate()
$wsRiders.Range("E9").Value = 1
$wsRiders.Range("E10").Value = 1
$wsRiders.Range("E11").Value = 1
[void]$wsRiders.Range("H11").Select()

# --- Drivers: becomes the active tab, with a new last-selected cell ---
$wsDrivers = $wb.Worksheets.Item("Drivers")
[void]$wsDrivers.Activate()
[void]$wsDrivers.Range("E32").Select()
